# Rename the "magnesium_mM" column header (cell F1) to "mg_mM"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F1").Value = "mg_mM"

# Update the active selection to the renamed header cell
$ws.Range("F1").Select()
